$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Update column F width (OOXML width 25 -> 24; ColumnWidth COM property is offset by ~0.83)
$ws.Range("F1").ColumnWidth = 23.17

# Row 3: update VENTA (D3); POR CUMPLIR (E3) and CUMPLIMIENTO (F3) recompute
$ws.Range("D3").Value = 1778.11
$ws.Range("E3").Value = $ws.Range("C3").Value2 - $ws.Range("D3").Value2
$ws.Range("F3").Value = $ws.Range("D3").Value2 / $ws.Range("C3").Value2

# Row 15: update VENTA (D15); POR CUMPLIR (E15) and CUMPLIMIENTO (F15) recompute
$ws.Range("D15").Value = 1753.96
$ws.Range("E15").Value = $ws.Range("C15").Value2 - $ws.Range("D15").Value2
$ws.Range("F15").Value = $ws.Range("D15").Value2 / $ws.Range("C15").Value2

# Row 19 (TOTAL): recompute VENTA as sum of D2:D18, then POR CUMPLIR and CUMPLIMIENTO
$ws.Range("D19").Value = $ws.Application.WorksheetFunction.Sum($ws.Range("D2:D18"))
$ws.Range("E19").Value = $ws.Range("C19").Value2 - $ws.Range("D19").Value2
$ws.Range("F19").Value = $ws.Range("D19").Value2 / $ws.Range("C19").Value2
